# Auto-generated script to apply numeric updates to the Mateus_Profits workbook
# per the target diff. Each worksheet tab (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# receives direct cell value updates; a few cells are cleared (removed) or newly added
# to match the diff precisely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 62535.25
$ws.Range("I6").Value = 83344.5
$ws.Range("K6").Value = 250033.5
$ws.Range("M6").Value = -249921.5
$ws.Range("H17").Value = 11112757
$ws.Range("J17").Value = 11112757
$ws.Range("L17").Value = 33338271
$ws.Range("N17").Value = -33338607
$ws.Range("H33").Value = 520.04346
$ws.Range("I33").Value = 385.55
$ws.Range("K33").Value = 385.55
$ws.Range("M33").Value = -156.55
$ws.Range("H64").Value = 10621.667
$ws.Range("J64").Value = 11027.083
$ws.Range("L64").Value = 11027.083
$ws.Range("N64").Value = -11523.083
$ws.Range("H67").Value = 10621.667
$ws.Range("J67").Value = 11027.083
$ws.Range("L67").Value = 11027.083
$ws.Range("N67").Value = -12743.083
$ws.Range("H106").Value = 5320.769
$ws.Range("I106").Value = 5080
$ws.Range("K106").Value = 5080
$ws.Range("M106").Value = -4449
$ws.Range("H125").Value = 14671.923
$ws.Range("I125").Value = 16650.223
$ws.Range("K125").Value = 149852.007
$ws.Range("M125").Value = -147392.007
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2019.8
$ws.Range("I26").Value = 2019.8
$ws.Range("K26").Value = 2019.8
$ws.Range("M26").Value = -1689.8
$ws.Range("H32").Value = 3786.5264
$ws.Range("I32").Value = 3221.9387
$ws.Range("K32").Value = 3221.9387
$ws.Range("M32").Value = -2934.9387
$ws.Range("H45").Value = 91105.87
$ws.Range("I45").Value = 136218.73
$ws.Range("J45").Value = 6519.25
$ws.Range("K45").Value = 136218.73
$ws.Range("L45").Value = 6519.25
$ws.Range("M45").Value = -135841.73
$ws.Range("N45").Value = -7273.25
$ws.Range("H63").Value = 4411
$ws.Range("I63").Value = 2797.8
$ws.Range("J63").Value = 8444
$ws.Range("K63").Value = 2797.8
$ws.Range("L63").Value = 8444
$ws.Range("M63").Value = -2111.8
$ws.Range("N63").Value = -9816
$ws.Range("H66").Value = 4411
$ws.Range("I66").Value = 2797.8
$ws.Range("J66").Value = 8444
$ws.Range("K66").Value = 13989
$ws.Range("L66").Value = 42220
$ws.Range("M66").Value = -10557
$ws.Range("N66").Value = -49084
$ws.Range("H102").Value = 3009.9565
$ws.Range("I102").Value = 2196.45
$ws.Range("K102").Value = 2196.45
$ws.Range("M102").Value = -574.4499999999998
$ws.Range("H110").Value = 8186
$ws.Range("I110").Value = 4958.4287
$ws.Range("K110").Value = 4958.4287
$ws.Range("M110").Value = -2913.4287
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 1981.4688
$ws.Range("I122").Value = 1768.091
$ws.Range("K122").Value = 5304.272999999999
$ws.Range("M122").Value = -2854.272999999999
$ws.Range("H132").Value = 3970.587
$ws.Range("I132").Value = 3148.2896
$ws.Range("K132").Value = 9444.8688
$ws.Range("M132").Value = -6914.8688

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3918.1177
$ws.Range("I105").Value = 3584.1667
$ws.Range("K105").Value = 3584.1667
$ws.Range("M105").Value = -1837.1667
$ws.Range("H132").Value = 15000
$ws.Range("J132").Value = 15000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -25120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 257.26666
$ws.Range("I7").Value = 95.44444
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 95.44444
$ws.Range("L7").Value = 500
$ws.Range("M7").Value = 17.55556
$ws.Range("N7").Value = -726
$ws.Range("H16").Value = 2212.5386
$ws.Range("I16").Value = 1703.8667
$ws.Range("K16").Value = 1703.8667
$ws.Range("M16").Value = -1416.8667
$ws.Range("H51").Value = 19173.637
$ws.Range("J51").Value = 19173.637
$ws.Range("L51").Value = 19173.637
$ws.Range("N51").Value = -20645.637
$ws.Range("H58").Value = 5605.75
$ws.Range("I58").Value = 2688.111
$ws.Range("K58").Value = 2688.111
$ws.Range("M58").Value = -2485.111
$ws.Range("H61").Value = 19173.637
$ws.Range("J61").Value = 19173.637
$ws.Range("L61").Value = 19173.637
$ws.Range("N61").Value = -19869.637
$ws.Range("H113").Value = 2212.5386
$ws.Range("I113").Value = 1703.8667
$ws.Range("K113").Value = 1703.8667
$ws.Range("M113").Value = 466.1333
$ws.Range("H136").Value = 5605.75
$ws.Range("I136").Value = 2688.111
$ws.Range("K136").Value = 8064.333
$ws.Range("M136").Value = -5514.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1146
$ws.Range("J38").Value = 1332.5
$ws.Range("L38").Value = 3997.5
$ws.Range("N38").Value = -4691.5
$ws.Range("H80").Value = 5628.3335
$ws.Range("I80").Value = 5944
$ws.Range("K80").Value = 17832
$ws.Range("M80").Value = -16896
$ws.Range("H83").Value = 5628.3335
$ws.Range("I83").Value = 5944
$ws.Range("K83").Value = 53496
$ws.Range("M83").Value = -48816
$ws.Range("H140").Value = 5021.5
$ws.Range("I140").Value = 8043
$ws.Range("K140").Value = 24129
$ws.Range("M140").Value = -18949

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 630.2353000000001
$ws.Range("I2").Value = 43.75
$ws.Range("J2").Value = 2037.8
$ws.Range("K2").Value = 43.75
$ws.Range("L2").Value = 2037.8
$ws.Range("M2").Value = 69.25
$ws.Range("N2").Value = -2263.8
$ws.Range("H70").Value = 13674
$ws.Range("I70").Value = 10354.5
$ws.Range("K70").Value = 10354.5
$ws.Range("M70").Value = -10084.5
$ws.Range("H73").Value = 13674
$ws.Range("I73").Value = 10354.5
$ws.Range("K73").Value = 10354.5
$ws.Range("M73").Value = -9418.5
$ws.Range("H80").Value = 4445.909
$ws.Range("J80").Value = 4036.5
$ws.Range("L80").Value = 4036.5
$ws.Range("N80").Value = -6032.5
$ws.Range("H83").Value = 4445.909
$ws.Range("J83").Value = 4036.5
$ws.Range("L83").Value = 20182.5
$ws.Range("N83").Value = -30166.5
$ws.Range("H132").Value = 3399.4
$ws.Range("I132").Value = 2999.25
$ws.Range("K132").Value = 8997.75
$ws.Range("M132").Value = -6467.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1822.0834
$ws.Range("I16").Value = 1896.909
$ws.Range("J16").Value = 999
$ws.Range("K16").Value = 1896.909
$ws.Range("L16").Value = 999
$ws.Range("M16").Value = -1726.909
$ws.Range("N16").Value = -1339
$ws.Range("H22").Value = 2665.3333
$ws.Range("I22").Value = 3498.5
$ws.Range("K22").Value = 3498.5
$ws.Range("M22").Value = -3203.5
$ws.Range("H27").Value = 2665.3333
$ws.Range("I27").Value = 3498.5
$ws.Range("K27").Value = 3498.5
$ws.Range("M27").Value = -3391.5
$ws.Range("H100").Value = 2944059.5
$ws.Range("I100").Value = 5002210
$ws.Range("J100").Value = 3844.2856
$ws.Range("K100").Value = 5002210
$ws.Range("L100").Value = 3844.2856
$ws.Range("M100").Value = -5001669
$ws.Range("N100").Value = -4926.2856
$ws.Range("H132").Value = 19218
$ws.Range("I132").Value = 51404
$ws.Range("K132").Value = 154212
$ws.Range("M132").Value = -151682
$ws.Range("H133").Value = 72980
$ws.Range("J133").Value = 72980
$ws.Range("L133").Value = 72980
$ws.Range("N133").Value = -78040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1163.1333
$ws.Range("I100").Value = 1205.8182
$ws.Range("J100").Value = 1045.75
$ws.Range("K100").Value = 2411.6364
$ws.Range("L100").Value = 2091.5
$ws.Range("M100").Value = -1870.6364
$ws.Range("N100").Value = -3173.5
$ws.Range("H109").Value = 69996.664
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 69996.664
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 69996.664
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -72770.664
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 3465.5957
$ws.Range("I136").Value = 2697.8572
$ws.Range("K136").Value = 8093.571599999999
$ws.Range("M136").Value = -5543.571599999999
